# Update "想去人数" (column F) figures on both the "展览" sheet and the
# aggregated "全部类型" sheet, matching the latest scrape output.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 1593
    5  = 610
    7  = 10
    8  = 11374
    15 = 12339
    16 = 13009
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
